$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need to be
# force-formatted as Text first, otherwise Excel auto-converts the
# assigned string into a real number (losing exact text representation,
# e.g. trailing zeros like "0.340" or "1.00").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.288.45"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.680.62"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "522.86"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "146.61"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "0.576"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "2.700.41"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "0.340"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "3.153.43"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "60.390.57"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "21.36"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.698.78"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "351.88"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "10.56"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("D24").Value = "63.24"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "0.0₃0817"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "6.78"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "19.16"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "1.60"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "147.26"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").Value = "4.30"
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("D36").Value = "1.26"
$ws.Range("E36").Value = "  +8.41%  "
$ws.Range("D37").Value = "0.954"
$ws.Range("E37").Value = "  -5.77%  "
$ws.Range("D38").Value = "0.879"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +6.82%  "
$ws.Range("D40").Value = "36.95"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "284.36"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "0.0992"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "20.06"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.612"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "2.133.75"
$ws.Range("E47").Value = "  +5.62%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0541"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.88"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "0.0236"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "19.44"
$ws.Range("E51").Value = "  +4.27%  "
